$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Data Types" -> "Variables and data types"; End Date 45652 -> 45656; Total Days 1 -> 5
$ws.Range("A3").Value = "Variables and data types"
$ws.Range("D3").Value = 45656
$ws.Range("E3").Value = 5

# Row 4: add file name "controlstructures.py"; Start/End Date move to 45657/45657; Total Days 3 -> 1
$ws.Range("B4").Value = "controlstructures.py"
$ws.Range("C4").Value = 45657
$ws.Range("D4").Value = 45657
$ws.Range("E4").Value = 1

# Row 5: Start Date 45657 -> 45658; Total Days 2 -> 1 (End Date 45658 unchanged)
$ws.Range("C5").Value = 45658
$ws.Range("E5").Value = 1

# Column width adjustments (A: 18.6640625 -> 21, B: 14.6640625 -> 17.6640625 with bestFit)
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668
$ws.Columns.Item(2).ColumnWidth = 16.830729166666668

# Move the active selection to I3
$ws.Range("I3").Select()
